# Update the "取得日時" (retrieved timestamp) column for the newly
# re-fetched rows to reflect the latest scrape time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-24 06:34:39"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
